$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (LeetCode problems 3354, 3370, 2464) ---

# Row 193: 3354 - Make Array Elements Equal to Zero
$ws.Cells.Item(193,1).Value = 3354
$ws.Cells.Item(193,2).Value = "Make Array Elements Equal to Zero"
$ws.Cells.Item(193,3).Value = "#array #prefix-sum #simulation "
$ws.Cells.Item(193,4).Value = "easy"
$ws.Cells.Item(193,5).Value = 1
$ws.Cells.Item(193,6).Value = 0
$ws.Cells.Item(193,7).Value = 7
$ws.Cells.Item(193,8).Value = 45958
$ws.Cells.Item(193,8).NumberFormat = "m/d/yy"
$ws.Cells.Item(193,9).Value = 45958
$ws.Cells.Item(193,9).NumberFormat = "m/d/yy"

# Row 194: 3370 - Smallest Number With All Set Bits
$ws.Cells.Item(194,1).Value = 3370
$ws.Cells.Item(194,2).Value = "Smallest Number With All Set Bits"
$ws.Cells.Item(194,3).Value = "#bit-minipulation "
$ws.Cells.Item(194,4).Value = "easy"
$ws.Cells.Item(194,5).Value = 1
$ws.Cells.Item(194,6).Value = 0
$ws.Cells.Item(194,7).Value = 4
$ws.Cells.Item(194,8).Value = 45959
$ws.Cells.Item(194,8).NumberFormat = "m/d/yy"
$ws.Cells.Item(194,9).Value = 45959
$ws.Cells.Item(194,9).NumberFormat = "m/d/yy"

# Row 195: 2464 - Minimum Subarrays in a Valid Split
$ws.Cells.Item(195,1).Value = 2464
$ws.Cells.Item(195,2).Value = "Minimum Subarrays in a Valid Split"
$ws.Cells.Item(195,3).Value = "#dynamic-programming #array "
$ws.Cells.Item(195,4).Value = "medium"
$ws.Cells.Item(195,5).Value = 0
$ws.Cells.Item(195,6).Value = 1
$ws.Cells.Item(195,7).Value = 20
$ws.Cells.Item(195,8).Value = 45959
$ws.Cells.Item(195,8).NumberFormat = "m/d/yy"
$ws.Cells.Item(195,9).Value = 45959
$ws.Cells.Item(195,9).NumberFormat = "m/d/yy"

# --- Row heights for wrapped text rows ---
$ws.Rows.Item(193).RowHeight = 34
$ws.Rows.Item(194).RowHeight = 34
$ws.Rows.Item(195).RowHeight = 51

# --- View state: scroll position / active selection following the new rows ---
$excel.ActiveWindow.ScrollRow = 191
$ws.Range("G195").Select()
